{"js": "// Adds a new group member (\"Muhammad Hanzla 24065593\") to the list of\n// authors directly below \"Muhammad Naveed Ashfaq 24065592\", matching the\n// formatting conventions already used in that paragraph block:\n//   - the existing \"Muhammad Naveed Ashfaq 24065592\" line gains a trailing\n//     comma (it is no longer the last name in the list) and loses the\n//     italic/yellow-highlight paragraph-mark formatting it used to carry;\n//   - a new line \"                     Muhammad Hanzla 24065593\" is added\n//     (styled like the \"Asad Ali\" line above it), carrying a zero-length\n//     \"_GoBack\" bookmark right after the id, like Word leaves behind after\n//     the most recent edit position;\n//   - a new, empty paragraph follows, inheriting the italic + yellow\n//     highlight paragraph-mark formatting that used to sit on the\n//     \"Ashfaq\" line (this is what's left behind once the paragraph mark\n//     moves down after inserting the new content).\n\nconst rFonts = '<w:rFonts w:ascii=\"Times New Roman\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>';\n\n// Locate the paragraph that currently reads \"...Muhammad Naveed Ashfaq 24065592\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(p => p.text.indexOf(\"Muhammad Naveed Ashfaq 24065592\") !== -1);\nif (!target) {\n  throw new Error('Could not find the \"Muhammad Naveed Ashfaq 24065592\" paragraph.');\n}\n\n// Build the three paragraphs that replace the single target paragraph:\n//  1) the Ashfaq line, split so \"Ashfaq\" is its own run (proofing-style\n//     run split) and a trailing comma appended; plain paragraph mark.\n//  2) the new Hanzla line with a \"_GoBack\" bookmark after the id.\n//  3) an empty paragraph carrying the italic + yellow highlight paragraph\n//     mark that used to belong to the Ashfaq line.\nconst replacementOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  '<w:p>' +\n    '<w:pPr><w:rPr>' + rFonts + '</w:rPr></w:pPr>' +\n    '<w:r><w:rPr>' + rFonts + '</w:rPr><w:tab/><w:t xml:space=\"preserve\">         Muhammad Naveed </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr>' + rFonts + '</w:rPr><w:t>Ashfaq</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr>' + rFonts + '</w:rPr><w:t xml:space=\"preserve\"> 24065592</w:t></w:r>' +\n    '<w:r><w:rPr>' + rFonts + '</w:rPr><w:t>,</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n    '<w:pPr><w:rPr>' + rFonts + '</w:rPr></w:pPr>' +\n    '<w:r><w:rPr>' + rFonts + '</w:rPr><w:t xml:space=\"preserve\">                     Muhammad </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr>' + rFonts + '</w:rPr><w:t>Hanzla</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr>' + rFonts + '</w:rPr><w:t xml:space=\"preserve\"> 24065593</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n  '<w:p>' +\n    '<w:pPr><w:rPr>' + rFonts + '<w:i/><w:highlight w:val=\"yellow\"/></w:rPr></w:pPr>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntarget.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Adds a new group member (\"Muhammad Hanzla 24065593\") to the list of\n# authors directly below \"Muhammad Naveed Ashfaq 24065592\", matching the\n# formatting conventions already used in that paragraph block:\n#   - the existing \"Muhammad Naveed Ashfaq 24065592\" line gains a trailing\n#     comma (it is no longer the last name in the list) and loses the\n#     italic/yellow-highlight paragraph-mark formatting it used to carry;\n#   - a new line \"                     Muhammad Hanzla 24065593\" is added\n#     (styled like the \"Asad Ali\" line above it), carrying a zero-length\n#     \"_GoBack\" bookmark right after the id, like Word leaves behind after\n#     the most recent edit position;\n#   - a new, empty paragraph follows, inheriting the italic + yellow\n#     highlight paragraph-mark formatting that used to sit on the\n#     \"Ashfaq\" line (this is what's left behind once the paragraph mark\n#     moves down after inserting the new content).\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that currently reads \"...Muhammad Naveed Ashfaq 24065592\".\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"Muhammad Naveed Ashfaq 24065592\")\nif (-not $found) {\n    throw 'Could not find the \"Muhammad Naveed Ashfaq 24065592\" paragraph.'\n}\n$targetPara = $findRange.Paragraphs(1)\n$targetRange = $targetPara.Range\n\n$rFonts = '<w:rFonts w:ascii=\"Times New Roman\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>'\n\n# Build the three paragraphs that replace the single target paragraph:\n#  1) the Ashfaq line, split so \"Ashfaq\" is its own run (proofing-style\n#     run split) and a trailing comma appended; plain paragraph mark.\n#  2) the new Hanzla line with a \"_GoBack\" bookmark after the id.\n#  3) an empty paragraph carrying the italic + yellow highlight paragraph\n#     mark that used to belong to the Ashfaq line.\n$newBody =\n    '<w:p><w:pPr><w:rPr>' + $rFonts + '</w:rPr></w:pPr>' +\n    '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:tab/><w:t xml:space=\"preserve\">         Muhammad Naveed </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t>Ashfaq</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t xml:space=\"preserve\"> 24065592</w:t></w:r>' +\n    '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t>,</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p><w:pPr><w:rPr>' + $rFonts + '</w:rPr></w:pPr>' +\n    '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t xml:space=\"preserve\">                     Muhammad </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t>Hanzla</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t xml:space=\"preserve\"> 24065593</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '</w:p>' +\n    '<w:p><w:pPr><w:rPr>' + $rFonts + '<w:i/><w:highlight w:val=\"yellow\"/></w:rPr></w:pPr></w:p>'\n\n$ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $newBody + '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n$targetRange.InsertXML($ooxml)\n"}
